# Add a new worksheet "2025-07-01" with the day's ranking data,
# inserted right after the last existing date sheet ("2025-06-30").
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2025-07-01"

# Header row
$newSheet.Cells.Item(1, 1).Value2 = "rank"
$newSheet.Cells.Item(1, 2).Value2 = "title"
$newSheet.Cells.Item(1, 3).Value2 = "author"
$newSheet.Cells.Item(1, 4).Value2 = "latest_episode"

# Match the bold / bordered / centered header style used by the other date sheets
$headerStyleSource = $lastSheet.Range("A1")
$headerStyleSource.Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Ranking data rows
$newSheet.Cells.Item(2, 1).Value2 = 1
$newSheet.Cells.Item(2, 2).Value2 = '蜘蛛ですが、なにか？'
$newSheet.Cells.Item(2, 3).Value2 = 'かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)'
$newSheet.Cells.Item(2, 4).Value2 = '第75話その2'
$newSheet.Cells.Item(3, 1).Value2 = 2
$newSheet.Cells.Item(3, 2).Value2 = '帰ってください！ 阿久津さん'
$newSheet.Cells.Item(3, 3).Value2 = '長岡太一(著者)'
$newSheet.Cells.Item(3, 4).Value2 = '第191話'
$newSheet.Cells.Item(4, 1).Value2 = 3
$newSheet.Cells.Item(4, 2).Value2 = '衛宮さんちの今日のごはん'
$newSheet.Cells.Item(4, 3).Value2 = 'TAa(漫画) 只野まこと(料理監修) ＴＹＰＥ－ＭＯＯＮ(原作)'
$newSheet.Cells.Item(4, 4).Value2 = '第74話'
$newSheet.Cells.Item(5, 1).Value2 = 4
$newSheet.Cells.Item(5, 2).Value2 = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$newSheet.Cells.Item(5, 3).Value2 = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$newSheet.Cells.Item(5, 4).Value2 = '第１８話①'
$newSheet.Cells.Item(6, 1).Value2 = 5
$newSheet.Cells.Item(6, 2).Value2 = '異世界居酒屋「のぶ」'
$newSheet.Cells.Item(6, 3).Value2 = '蝉川夏哉(原作) ヴァージニア二等兵(漫画) 転(キャラクター原案)'
$newSheet.Cells.Item(6, 4).Value2 = '第122話'
$newSheet.Cells.Item(7, 1).Value2 = 6
$newSheet.Cells.Item(7, 2).Value2 = '異世界建国記'
$newSheet.Cells.Item(7, 3).Value2 = 'ＫＯＩＺＵＭＩ(漫画) 桜木桜(原作) 屡那(キャラクター原案)'
$newSheet.Cells.Item(7, 4).Value2 = '第76話'
$newSheet.Cells.Item(8, 1).Value2 = 7
$newSheet.Cells.Item(8, 2).Value2 = '一億年ボタンを連打した俺は、気付いたら最強になっていた ～落第剣士の学院無双～'
$newSheet.Cells.Item(8, 3).Value2 = '士土幽太郎(漫画) 月島秀一(原作) もきゅ(キャラクター原案)'
$newSheet.Cells.Item(8, 4).Value2 = '第40話-1'
$newSheet.Cells.Item(9, 1).Value2 = 8
$newSheet.Cells.Item(9, 2).Value2 = '尾守つみきと奇日常。'
$newSheet.Cells.Item(9, 3).Value2 = '森下みゆ'
$newSheet.Cells.Item(9, 4).Value2 = '第57話 つみきさん達とパーティー。'
$newSheet.Cells.Item(10, 1).Value2 = 9
$newSheet.Cells.Item(10, 2).Value2 = 'クラス最安値で売られた俺は、実は最強パラメーター'
$newSheet.Cells.Item(10, 3).Value2 = 'カンブリア爆発太郎(漫画) RYOMA(原作) 黒井ススム(キャラクター原案)'
$newSheet.Cells.Item(10, 4).Value2 = '第35話-2'
$newSheet.Cells.Item(11, 1).Value2 = 10
$newSheet.Cells.Item(11, 2).Value2 = '狂戦士なモブ、無自覚に本編を破壊する'
$newSheet.Cells.Item(11, 3).Value2 = '漫画：佐藤良亮 原作：なるのるな キャラクター原案：霜月えいと'
$newSheet.Cells.Item(11, 4).Value2 = '第11話 ①'
$newSheet.Cells.Item(12, 1).Value2 = 11
$newSheet.Cells.Item(12, 2).Value2 = 'レベル１から始まる召喚無双'
$newSheet.Cells.Item(12, 3).Value2 = '漫画：七桃りお 原作：白石新 キャラクター原案：夕薙'
$newSheet.Cells.Item(12, 4).Value2 = '第32話前半'
$newSheet.Cells.Item(13, 1).Value2 = 12
$newSheet.Cells.Item(13, 2).Value2 = '不遇皇子は天才錬金術師～皇帝なんて柄じゃないので弟妹を可愛がりたい～@COMIC'
$newSheet.Cells.Item(13, 3).Value2 = '長先ザワ（漫画） うめー（原作） 瑛来イチ（構成） 雨銛（構成） かわく（キャラクター原案）'
$newSheet.Cells.Item(13, 4).Value2 = '第8話 ①'
$newSheet.Cells.Item(14, 1).Value2 = 13
$newSheet.Cells.Item(14, 2).Value2 = '異世界おじさん'
$newSheet.Cells.Item(14, 3).Value2 = '殆ど死んでいる(著者)'
$newSheet.Cells.Item(14, 4).Value2 = '第69話'
$newSheet.Cells.Item(15, 1).Value2 = 14
$newSheet.Cells.Item(15, 2).Value2 = 'ひきこもりの俺がかわいいギルドマスターに世話を焼かれまくったって別にいいだろう?'
$newSheet.Cells.Item(15, 3).Value2 = '漫画：ミト 原作：東條功一 イラスト：にもし'
$newSheet.Cells.Item(15, 4).Value2 = '第17話'
$newSheet.Cells.Item(16, 1).Value2 = 15
$newSheet.Cells.Item(16, 2).Value2 = 'S級パーティーから追放された狩人、実は世界最強 ～射程9999の男、帝国の狙撃手として無双する～'
$newSheet.Cells.Item(16, 3).Value2 = '漫画：カズミヤアキラ 原作：茨木野 キャラクター原案：へいろー'
$newSheet.Cells.Item(16, 4).Value2 = '第8話 ③'
$newSheet.Cells.Item(17, 1).Value2 = 16
$newSheet.Cells.Item(17, 2).Value2 = '元最強の剣士は、異世界魔法に憧れる'
$newSheet.Cells.Item(17, 3).Value2 = '漫画：天乃ちはる 原作：紅月シン キャラクター原案：necömi'
$newSheet.Cells.Item(17, 4).Value2 = '第72話前半'
$newSheet.Cells.Item(18, 1).Value2 = 17
$newSheet.Cells.Item(18, 2).Value2 = 'アラフォーになった最強の英雄たち、再び戦場で無双する!!'
$newSheet.Cells.Item(18, 3).Value2 = '漫画：戸田泰成 原作：岸馬きらく 構成協力：森小太郎 キャラクター原案：peroshi'
$newSheet.Cells.Item(18, 4).Value2 = '第23話 前編'
$newSheet.Cells.Item(19, 1).Value2 = 18
$newSheet.Cells.Item(19, 2).Value2 = 'スライムは最強たる可能性を秘めている～２回目の人生、ちゃんとスライムと向き合います～@COMIC'
$newSheet.Cells.Item(19, 3).Value2 = 'オサフネオウジ（漫画） 犬型大（原作） 風花風花（キャラクター原案）'
$newSheet.Cells.Item(19, 4).Value2 = '第6話'
$newSheet.Cells.Item(20, 1).Value2 = 19
$newSheet.Cells.Item(20, 2).Value2 = '能あるオーガは角を隠す'
$newSheet.Cells.Item(20, 3).Value2 = '漫画家： 蒼葉 結 原作： 津野瀬 文'
$newSheet.Cells.Item(20, 4).Value2 = '第9話 前編'
$newSheet.Cells.Item(21, 1).Value2 = 20
$newSheet.Cells.Item(21, 2).Value2 = 'クセ強彼女は床にいざなう'
$newSheet.Cells.Item(21, 3).Value2 = '須河篤志(著者)'
$newSheet.Cells.Item(21, 4).Value2 = '第13話前半'
$newSheet.Cells.Item(22, 1).Value2 = 21
$newSheet.Cells.Item(22, 2).Value2 = 'ニチアサ好きのオタクが悪役生徒に転生した結果、破滅フラグが崩壊していく件について'
$newSheet.Cells.Item(22, 3).Value2 = '烏丸英（原作） どんぐりす（漫画）'
$newSheet.Cells.Item(22, 4).Value2 = '第14話（後編）急襲…事件の始まり'
$newSheet.Cells.Item(23, 1).Value2 = 22
$newSheet.Cells.Item(23, 2).Value2 = '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～'
$newSheet.Cells.Item(23, 3).Value2 = '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)'
$newSheet.Cells.Item(23, 4).Value2 = '第8話-2'
$newSheet.Cells.Item(24, 1).Value2 = 23
$newSheet.Cells.Item(24, 2).Value2 = 'ちはるくんは女装をしたくない！'
$newSheet.Cells.Item(24, 3).Value2 = '翁丸ジョン'
$newSheet.Cells.Item(24, 4).Value2 = '【第19話】男装宗と交流したくない！その一'
$newSheet.Cells.Item(25, 1).Value2 = 24
$newSheet.Cells.Item(25, 2).Value2 = 'ぽんドロイド！ はまさん'
$newSheet.Cells.Item(25, 3).Value2 = 'はれやまはれぞう(著者)'
$newSheet.Cells.Item(25, 4).Value2 = '第2話'
$newSheet.Cells.Item(26, 1).Value2 = 25
$newSheet.Cells.Item(26, 2).Value2 = '大魔導士と呼ばれた侯爵令嬢〜世界が汚いので掃除していただけなんですけど……〜@COMIC'
$newSheet.Cells.Item(26, 3).Value2 = 'さいピン（漫画） K1you（原作） パルプピロシ（キャラクター原案）'
$newSheet.Cells.Item(26, 4).Value2 = '第6話'
$newSheet.Cells.Item(27, 1).Value2 = 26
$newSheet.Cells.Item(27, 2).Value2 = 'クソ女に幸あれ'
$newSheet.Cells.Item(27, 3).Value2 = '岸川瑞樹'
$newSheet.Cells.Item(27, 4).Value2 = '第58話'
$newSheet.Cells.Item(28, 1).Value2 = 27
$newSheet.Cells.Item(28, 2).Value2 = '魔術師クノンは見えている'
$newSheet.Cells.Item(28, 3).Value2 = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$newSheet.Cells.Item(28, 4).Value2 = '第38話①'
$newSheet.Cells.Item(29, 1).Value2 = 28
$newSheet.Cells.Item(29, 2).Value2 = 'リビルドワールド'
$newSheet.Cells.Item(29, 3).Value2 = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$newSheet.Cells.Item(29, 4).Value2 = '第70話③'
$newSheet.Cells.Item(30, 1).Value2 = 29
$newSheet.Cells.Item(30, 2).Value2 = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$newSheet.Cells.Item(30, 3).Value2 = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$newSheet.Cells.Item(30, 4).Value2 = '第16話後半：「ストリア見聞録」'
$newSheet.Cells.Item(31, 1).Value2 = 30
$newSheet.Cells.Item(31, 2).Value2 = 'ワンパンマン'
$newSheet.Cells.Item(31, 3).Value2 = '原作/ＯＮＥ 作画/村田雄介'
$newSheet.Cells.Item(31, 4).Value2 = '201撃目'
$newSheet.Cells.Item(32, 1).Value2 = 31
$newSheet.Cells.Item(32, 2).Value2 = '千年英雄'
$newSheet.Cells.Item(32, 3).Value2 = '原作/福島航平 作画/中村ゆきひろ'
$newSheet.Cells.Item(32, 4).Value2 = '12話②'
$newSheet.Cells.Item(33, 1).Value2 = 32
$newSheet.Cells.Item(33, 2).Value2 = 'サーシャちゃんとクラスメイトオタクくん'
$newSheet.Cells.Item(33, 3).Value2 = 'はぐはぐ(著者)'
$newSheet.Cells.Item(33, 4).Value2 = '第80話'
$newSheet.Cells.Item(34, 1).Value2 = 33
$newSheet.Cells.Item(34, 2).Value2 = '冒険者ギルドが十二歳からしか入れなかったので、サバよみました。'
$newSheet.Cells.Item(34, 3).Value2 = 'GUNP（漫画） KAME （原作） ox （キャラクター原案）'
$newSheet.Cells.Item(34, 4).Value2 = '第12話前半'
$newSheet.Cells.Item(35, 1).Value2 = 34
$newSheet.Cells.Item(35, 2).Value2 = '《魔力無限》のマナポーター ～パーティの魔力を全て供給していたのに、勇者に追放されました。魔力不足で聖剣が使えないと焦っても、メンバー全員が勇者を見限ったのでもう遅い～'
$newSheet.Cells.Item(35, 3).Value2 = '漫画：伊藤ひずみ 原作：アトハ キャラクター原案：夕薙'
$newSheet.Cells.Item(35, 4).Value2 = '第10話 ②'
$newSheet.Cells.Item(36, 1).Value2 = 35
$newSheet.Cells.Item(36, 2).Value2 = '淫獄団地'
$newSheet.Cells.Item(36, 3).Value2 = '搾精研究所(原作) 丈山雄為(漫画)'
$newSheet.Cells.Item(36, 4).Value2 = '第48話（前編）'
$newSheet.Cells.Item(37, 1).Value2 = 36
$newSheet.Cells.Item(37, 2).Value2 = 'ダンジョンの幼なじみ'
$newSheet.Cells.Item(37, 3).Value2 = '久真やすひさ(著者)'
$newSheet.Cells.Item(37, 4).Value2 = '第54話'
$newSheet.Cells.Item(38, 1).Value2 = 37
$newSheet.Cells.Item(38, 2).Value2 = 'オトナを知りたい女友達'
$newSheet.Cells.Item(38, 3).Value2 = '望公太(原作) ぷよちゃ(作画)'
$newSheet.Cells.Item(38, 4).Value2 = '第7話 前半'
$newSheet.Cells.Item(39, 1).Value2 = 38
$newSheet.Cells.Item(39, 2).Value2 = '田舎の黒ギャルJKと結婚しました'
$newSheet.Cells.Item(39, 3).Value2 = 'カヅチ(著者)'
$newSheet.Cells.Item(39, 4).Value2 = '第16話後半'
$newSheet.Cells.Item(40, 1).Value2 = 39
$newSheet.Cells.Item(40, 2).Value2 = '異世界迷宮のオーパーツ'
$newSheet.Cells.Item(40, 3).Value2 = '三狛ハル(著者)'
$newSheet.Cells.Item(40, 4).Value2 = '第1話：先端にお椀がついた棒'
$newSheet.Cells.Item(41, 1).Value2 = 40
$newSheet.Cells.Item(41, 2).Value2 = 'ダークサモナーとデキている'
$newSheet.Cells.Item(41, 3).Value2 = '車王(著者)'
$newSheet.Cells.Item(41, 4).Value2 = '第72話'
$newSheet.Cells.Item(42, 1).Value2 = 41
$newSheet.Cells.Item(42, 2).Value2 = 'ギャルゲーマーに褒められたい'
$newSheet.Cells.Item(42, 3).Value2 = 'げしゅまろ(著者)'
$newSheet.Cells.Item(42, 4).Value2 = '34話'
$newSheet.Cells.Item(43, 1).Value2 = 42
$newSheet.Cells.Item(43, 2).Value2 = '転生貴族の異世界冒険録 ～カインのやりすぎギルド日記～'
$newSheet.Cells.Item(43, 3).Value2 = '原作：夜州 漫画：佐々木あかね・香本セトラ キャラクター原案：藻'
$newSheet.Cells.Item(43, 4).Value2 = '第47話'
$newSheet.Cells.Item(44, 1).Value2 = 43
$newSheet.Cells.Item(44, 2).Value2 = '異世界のんびり農家'
$newSheet.Cells.Item(44, 3).Value2 = '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)'
$newSheet.Cells.Item(44, 4).Value2 = '第301話'
$newSheet.Cells.Item(45, 1).Value2 = 44
$newSheet.Cells.Item(45, 2).Value2 = '小さめの魔法師匠と大きめの魔法少女。report：3'
$newSheet.Cells.Item(45, 3).Value2 = 'とりから'
$newSheet.Cells.Item(45, 4).Value2 = '第36話の3'
$newSheet.Cells.Item(46, 1).Value2 = 45
$newSheet.Cells.Item(46, 2).Value2 = '男嫌いな美人姉妹を名前も告げずに助けたら一体どうなる?'
$newSheet.Cells.Item(46, 3).Value2 = 'みょん(原作) 司馬淳子(漫画) ぎうにう(キャラクターデザイン)'
$newSheet.Cells.Item(46, 4).Value2 = '第22話'
$newSheet.Cells.Item(47, 1).Value2 = 46
$newSheet.Cells.Item(47, 2).Value2 = '聖液鍛冶屋のエロランタ'
$newSheet.Cells.Item(47, 3).Value2 = 'しげきっくす(著者)'
$newSheet.Cells.Item(47, 4).Value2 = '第6話後半'
$newSheet.Cells.Item(48, 1).Value2 = 47
$newSheet.Cells.Item(48, 2).Value2 = 'ロメリア戦記～伯爵令嬢、魔王を倒した後も人類やばそうだから軍隊組織する～'
$newSheet.Cells.Item(48, 3).Value2 = '漫画：上戸 亮 原作：有山リョウ(小学館「ガガガブックス」刊) キャラクター原案：コダマ'
$newSheet.Cells.Item(48, 4).Value2 = '第13話「助けてくれる人々」①'
$newSheet.Cells.Item(49, 1).Value2 = 48
$newSheet.Cells.Item(49, 2).Value2 = '老後に備えて異世界で８万枚の金貨を貯めます'
$newSheet.Cells.Item(49, 3).Value2 = 'FUNA 東西 モトエ恵介'
$newSheet.Cells.Item(49, 4).Value2 = '第118話　会談［その3］'
$newSheet.Cells.Item(50, 1).Value2 = 49
$newSheet.Cells.Item(50, 2).Value2 = 'まんきつしたい常連さん'
$newSheet.Cells.Item(50, 3).Value2 = 'しんみりん(著者)'
$newSheet.Cells.Item(50, 4).Value2 = '第45話前編'
$newSheet.Cells.Item(51, 1).Value2 = 50
$newSheet.Cells.Item(51, 2).Value2 = 'はずれスキル念動力（ただしレベルＭＡＸ）で無双する～手をかざすだけです。詠唱とか必殺技とかいりません。念じるだけで倒せます～'
$newSheet.Cells.Item(51, 3).Value2 = 'さとう うなぽっぽ トダフミト'
$newSheet.Cells.Item(51, 4).Value2 = '9話①'

Write-Output "Added sheet $($newSheet.Name) with $($newSheet.UsedRange.Rows.Count) rows"
